$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.336.32'
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").Value = '2.481.89'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'525.30"
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("D6").Value = "'133.48"
$ws.Range("E6").Value = '  -3.64%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("D9").Value = "'0.1000"
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("D11").Value = "'5.42"
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("D13").Value = '2.922.92'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("D14").Value = '58.313.46'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("D15").Value = "'22.39"
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("E16").Value = '  -2.06%  '
$ws.Range("D17").Value = '2.486.02'
$ws.Range("D18").Value = "'10.91"
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Value = "'4.19"
$ws.Range("E19").Value = '  -2.70%  '
$ws.Range("D20").Value = "'321.40"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = "'5.80"
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").Value = "'64.33"
$ws.Range("E23").Value = '  -1.73%  '
$ws.Range("D24").Value = "'0.412"
$ws.Range("E24").Value = '  -3.24%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("D27").Value = "'7.47"
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("D28").Value = '0.0₃0751'
$ws.Range("E28").Value = '  -3.87%  '
$ws.Range("D29").Value = "'6.41"
$ws.Range("E29").Value = '  -5.31%  '
$ws.Range("E30").Value = '  -4.67%  '
$ws.Range("D31").Value = "'166.35"
$ws.Range("E31").Value = '  -1.40%  '
$ws.Range("E32").Value = '  -4.81%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("E36").Value = '  -9.18%  '
$ws.Range("E37").Value = '  -4.09%  '
$ws.Range("E38").Value = '  -4.22%  '
$ws.Range("D39").Value = "'0.797"
$ws.Range("E39").Value = '  -3.39%  '
$ws.Range("E40").Value = '  -3.68%  '
$ws.Range("D41").Value = "'277.85"
$ws.Range("E41").Value = '  -2.62%  '
$ws.Range("D42").Value = "'4.96"
$ws.Range("E42").Value = '  -5.39%  '
$ws.Range("D43").Value = "'0.594"
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").Value = "'127.50"
$ws.Range("E44").Value = '  -3.29%  '
$ws.Range("D45").Value = "'0.0912"
$ws.Range("E45").Value = '  -2.41%  '
$ws.Range("D46").Value = "'0.0495"
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("E47").Value = '  -2.62%  '
$ws.Range("D48").Value = "'17.21"
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("D49").Value = '1.741.95'
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").Value = "'0.972"
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").Value = "'4.67"
$ws.Range("E51").Value = '  -2.11%  '
